# Insert a new data row before current row 56, shifting rows 56..88 down to 57..89.
# The new row 56 keeps the same Volumen/Precio columns (J,K,L,M,P) that the old
# row 56 had, but gets a new date (D56 = 44488). All other columns (A,B,C,E,F,
# G,H,I,N,O,Q,R) are identical across every data row in this sheet, so the newly
# inserted row is simply filled with those same constant values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value = 8
$ws.Cells.Item(56, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(56, 3).Value = "Coquimbo"
$ws.Cells.Item(56, 4).Value = 44488
$ws.Cells.Item(56, 5).Value = 4
$ws.Cells.Item(56, 6).Value = 100112040
$ws.Cells.Item(56, 7).Value = "Cilantro"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 3200
$ws.Cells.Item(56, 11).Value = 1300
$ws.Cells.Item(56, 12).Value = 1500
$ws.Cells.Item(56, 13).Value = 1400
$ws.Cells.Item(56, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(56, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(56, 16).Value = 933
$ws.Cells.Item(56, 17).Value = 1.5
$ws.Cells.Item(56, 18).Value = "Hortaliza"

$ws.Cells.Item(56, 4).NumberFormat = $ws.Cells.Item(57, 4).NumberFormat
